# Ajout note sur word
# Insert 4 new paragraphs right after the "Pour le bit de signe..." paragraph
# (and right before the pre-existing trailing empty paragraph):
#   1) empty paragraph
#   2) empty paragraph
#   3) bold, 16pt: "Il faut faire un ReLu sur le résultat global du multiply and accumulate."
#   4) bold, 16pt: "Il y a une multiplication de l'entrée 1 de la matrice de l'image avec
#      le poids 1 de la mémoire et ainsi de suite jusqu'au poids 784 fois l'entrée 784…"

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document body (the empty
# paragraph that immediately follows "Pour le bit de signe..."). The new
# content must be inserted right before it so it stays last.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
$insertionPoint = $d.Range($trailing.Range.Start, $trailing.Range.Start)

# Push in the four paragraphs (as plain, unformatted text first) in a single
# edit so Word lays out the paragraph marks/tab stops consistently with the
# surrounding "Normal" text.
$insertionPoint.InsertBefore("`r`rIl faut faire un ReLu sur le résultat global du multiply and accumulate.`rIl y a une multiplication de l" + [char]0x2019 + "entrée 1 de la matrice de l" + [char]0x2019 + "image avec le poids 1 de la mémoire et ainsi de suite jusqu" + [char]0x2019 + "au poids 784 fois l" + [char]0x2019 + "entrée 784" + [char]0x2026 + "`r")

# The trailing (pre-existing) empty paragraph has shifted down by 4; re-grab
# everything by index so the bold run can be applied to the right two
# paragraphs.
$idxTrailing = $d.Paragraphs.Count
$idxReluPara = $idxTrailing - 2
$idxMacPara  = $idxTrailing - 1

$reluPara = $d.Paragraphs($idxReluPara)
$reluPara.Range.Font.Bold = $true
$reluPara.Range.Font.BoldBi = $true
$reluPara.Range.Font.Size = 16
$reluPara.Range.Font.SizeBi = 16

$macPara = $d.Paragraphs($idxMacPara)
$macPara.Range.Font.Bold = $true
$macPara.Range.Font.BoldBi = $true
$macPara.Range.Font.Size = 16
$macPara.Range.Font.SizeBi = 16
